$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 11:06"

# Rows that changed between the two data snapshots: either the country's
# stats were refreshed, or countries swapped places in the (cases-sorted)
# ranking and therefore the country name + stats at that row changed.
# Columns: RowNum, Pais, Casos totales, Nuevos casos, Casos activos,
#          Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(6, "Rusia", 627646, 6852, 393352, 225325, 0, 188, 8969),
    @(7, "India", 509753, 307, 296030, 198023, 0, 11, 15700),
    @(17, "Alemania", 194399, 0, 177500, 7873, 0, 0, 9026),
    @(20, "Banglades", 133978, 3504, 54318, 77965, 0, 34, 1695),
    @(32, "Indonesia", 52812, 1385, 21909, 28183, 0, 37, 2720),
    @(41, "Filipinas", 34803, 730, 9430, 24137, 0, 12, 1236),
    @(42, "Polonia", 33714, 319, 19972, 12307, 0, 6, 1435),
    @(54, "Kazajistan", 20319, 569, 12589, 7580, 0, 0, 150),
    @(69, "Chequia", 11044, 6, 7669, 3026, 0, 0, 349),
    @(81, "El Salvador", 5727, 210, 3447, 2137, 0, 10, 143),
    @(82, "Haiti", 5722, 179, 641, 4983, 0, 2, 98),
    @(83, "Kenia", 5533, 0, 1905, 3491, 0, 0, 137),
    @(114, "Lituania", 1813, 5, 1503, 232, 0, 0, 78),
    @(117, "Eslovaquia", 1657, 14, 1455, 174, 0, 0, 28),
    @(118, "Estado de Palestina", 1624, 67, 446, 1175, 0, 0, 3),
    @(119, "Guinea-Bisau", 1614, 0, 191, 1401, 0, 0, 22),
    @(120, "Eslovenia", 1572, 14, 1376, 87, 0, 0, 109),
    @(130, "Niger", 1062, 3, 924, 71, 0, 0, 67),
    @(141, "Uganda", 848, 15, 761, 87, 0, 0, 0),
    @(201, "Laos", 19, 0, 19, 0, 0, 0, 0),
    @(202, "Santa Lucia", 19, 0, 19, 0, 0, 0, 0),
    @(203, "Dominica", 18, 0, 18, 0, 0, 0, 0),
    @(204, "Fiyi", 18, 0, 18, 0, 0, 0, 0),
    @(208, "Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
    @(209, "Groenlandia", 13, 0, 13, 0, 0, 0, 0),
    @(212, "Montserrat", 11, 0, 10, 0, 0, 0, 1),
    @(213, "Seychelles", 11, 0, 11, 0, 0, 0, 0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}

Write-Output ("Updated " + $rows.Count + " rows")
